$d = $word.ActiveDocument

$tbl = $d.Tables.Item(1)

# Row 5 (version 2.8.1): ณัฐนันท์ (QA) -> ณัฐดนัย (DM); กิตติพศ (SP) -> วิรัตน์ (TL)
$cell1 = $tbl.Cell(5, 4)
$cell1.Range.Find.Execute("ณัฐนันท์", $true, $false, $false, $false, $false, $true, 0, $false, "ณัฐดนัย", 1)
$cell1.Range.Find.Execute(" (QA)", $true, $false, $false, $false, $false, $true, 0, $false, " (DM)", 1)

$cell2 = $tbl.Cell(5, 5)
$cell2.Range.Find.Execute("กิตติพศ ", $true, $false, $false, $false, $false, $true, 0, $false, "วิรัตน์", 1)
$cell2.Range.Find.Execute("(SP)", $true, $false, $false, $false, $false, $true, 0, $false, " (TL)", 1)

# Row 6 (version 2.4.1): ณัฐนันท์ (QA) -> วิรัตน์ (TL); กิตติพศ (SP) -> วริศรา (D)
$cell3 = $tbl.Cell(6, 4)
$cell3.Range.Find.Execute("ณัฐนันท์", $true, $false, $false, $false, $false, $true, 0, $false, "วิรัตน์", 1)
$cell3.Range.Find.Execute(" (QA)", $true, $false, $false, $false, $false, $true, 0, $false, " (TL)", 1)

$cell4 = $tbl.Cell(6, 5)
$cell4.Range.Find.Execute("กิตติพศ ", $true, $false, $false, $false, $false, $true, 0, $false, "วริศรา", 1)
$cell4.Range.Find.Execute("(SP)", $true, $false, $false, $false, $false, $true, 0, $false, " (D)", 1)

Write-Output "Done"
